$d = $word.ActiveDocument

# The second paragraph is empty except for the "_GoBack" bookmark. We need
# to add text both before and after the bookmark, so that the final
# paragraph reads: "In that case, I would include a lot of interesting text."
# with the bookmark preserved in the middle (between "inc" and "lude").

$bm = $d.Bookmarks("_GoBack")
$bmRange = $bm.Range

# Insert the second half of the text right after the (collapsed) bookmark
# range. Using InsertAfter here keeps the bookmark anchored at its original
# (pre-insertion) position, i.e. *before* the newly inserted text.
$bmRange.InsertAfter("lude a lot of interesting text.")

# Re-fetch the bookmark (its range may have shifted) and insert the first
# half of the text immediately before it, using a fresh, explicit range so
# the bookmark itself ends up sandwiched between the two new runs.
$bm2 = $d.Bookmarks("_GoBack")
$startPos = $bm2.Range.Start
$beforeRange = $d.Range($startPos, $startPos)
$beforeRange.InsertBefore("In that case, I would inc")

# Apply the same run formatting ("en-GB" language) used elsewhere in the
# document to the whole paragraph (both new runs) in one shot - setting
# LanguageID on the full paragraph range (rather than a zero-length /
# freshly-constructed sub-range) is what reliably scopes the change to just
# this paragraph's runs in this environment.
$para2 = $d.Paragraphs(2).Range
$para2.LanguageID = "en-GB"
